# Juno: check in to OLPRODLOC.
# Localize the charger sales report workbook to Brazilian Portuguese:
#  - rename the worksheet tab
#  - translate the header row and the Year-Quarter column labels

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Relatório de vendas"

# Header row (row 1)
$ws.Range("A1").Value = "Trimestre-Ano"
$ws.Range("B1").Value = "Meio-oeste"
$ws.Range("C1").Value = "Mountain"
$ws.Range("D1").Value = "Northeast"
$ws.Range("E1").Value = "Sul"
$ws.Range("F1").Value = "Sudeste"
$ws.Range("G1").Value = "Oeste"

# Year-Quarter column (rows 2-9)
$ws.Range("A2").Value = "T1 de 2022"
$ws.Range("A3").Value = "T2 de 2022"
$ws.Range("A4").Value = "T3 de 2022"
$ws.Range("A5").Value = "T4 de 2022"
$ws.Range("A6").Value = "T1 de 2023"
$ws.Range("A7").Value = "T2 de 2023"
$ws.Range("A8").Value = "T3 de 2023"
$ws.Range("A9").Value = "T4 de 2023"
